$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the value to be stored as literal text, even if it
    # looks like a number (e.g. "0.999"), by assigning a text
    # formula and then converting it to a static value via
    # copy / paste-special (values only). This avoids Excel
    # auto-converting numeric-looking strings to real numbers.
    $escaped = $val -replace '"', '""'
    $ws.Range($addr).Formula = '="' + $escaped + '"'
    $ws.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

$ws.Range("D2").Value = "62.133.97"
$ws.Range("E2").Value = "  -2.12%  "
$ws.Range("D3").Value = "3.423.79"
$ws.Range("E3").Value = "  -1.71%  "
$ws.Range("E4").Value = "  +0.02%  "
Set-TextValue "D5" "578.48"
$ws.Range("E5").Value = "  -0.68%  "
Set-TextValue "D6" "152.44"
$ws.Range("E6").Value = "  +3.22%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +1.19%  "
$ws.Range("E9").Value = "  +4.95%  "
Set-TextValue "D10" "0.124"
$ws.Range("E10").Value = "  -0.85%  "
Set-TextValue "D11" "0.417"
$ws.Range("E11").Value = "  +2.95%  "
$ws.Range("D12").Value = "4.011.64"
$ws.Range("E12").Value = "  -1.63%  "
$ws.Range("E13").Value = "  +0.56%  "
Set-TextValue "D14" "28.73"
$ws.Range("E14").Value = "  -3.28%  "
$ws.Range("D15").Value = "3.427.70"
$ws.Range("E15").Value = "  -1.60%  "
$ws.Range("E16").Value = "  -0.57%  "
$ws.Range("D17").Value = "62.182.00"
$ws.Range("E17").Value = "  -2.09%  "
Set-TextValue "D18" "6.52"
$ws.Range("E18").Value = "  +1.86%  "
$ws.Range("E19").Value = "  +0.58%  "
Set-TextValue "D20" "8.94"
$ws.Range("E20").Value = "  -4.35%  "
Set-TextValue "D21" "382.92"
$ws.Range("E21").Value = "  -2.06%  "
Set-TextValue "D22" "0.571"
$ws.Range("E22").Value = "  +0.62%  "
Set-TextValue "D23" "75.16"
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").Value = "3.562.80"
$ws.Range("E25").Value = "  -1.61%  "
Set-TextValue "D26" "0.0000112"
$ws.Range("E26").Value = "  -3.93%  "
$ws.Range("E27").Value = "  -0.14%  "
Set-TextValue "D28" "7.69"
$ws.Range("E28").Value = "  +0.42%  "
Set-TextValue "D29" "0.999"
$ws.Range("E29").Value = "  -0.15%  "
Set-TextValue "D30" "7.90"
$ws.Range("E30").Value = "  -4.44%  "
Set-TextValue "D31" "2.12"
$ws.Range("E31").Value = "  -1.20%  "
Set-TextValue "D32" "0.999"
$ws.Range("E32").Value = "  -0.11%  "
Set-TextValue "D33" "1.34"
$ws.Range("E33").Value = "  -1.79%  "
Set-TextValue "D34" "23.22"
$ws.Range("E34").Value = "  -1.50%  "
Set-TextValue "D35" "5.46"
$ws.Range("E35").Value = "  +2.13%  "
$ws.Range("E36").Value = "  +0.13%  "
Set-TextValue "D37" "6.93"
$ws.Range("E37").Value = "  -3.39%  "
$ws.Range("B38").Value = "EnergySwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D38" "31.11"
$ws.Range("E38").Value = "  -2.17%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D39" "168.57"
$ws.Range("E39").Value = "  -0.65%  "
$ws.Range("D40").Value = "3.456.37"
$ws.Range("E40").Value = "  -1.86%  "
Set-TextValue "D41" "0.0784"
$ws.Range("E41").Value = "  +2.38%  "
Set-TextValue "D42" "42.78"
$ws.Range("E42").Value = "  +0.88%  "
Set-TextValue "D43" "0.778"
$ws.Range("E43").Value = "  -2.89%  "
$ws.Range("E44").Value = "  -0.26%  "
Set-TextValue "D45" "1.68"
$ws.Range("E45").Value = "  -3.36%  "
$ws.Range("E46").Value = "  -3.13%  "
$ws.Range("D47").Value = "2.542.97"
$ws.Range("E47").Value = "  -2.65%  "
Set-TextValue "D48" "6.87"
$ws.Range("E48").Value = "  +0.66%  "
Set-TextValue "D49" "2.20"
$ws.Range("E49").Value = "  -4.85%  "
Set-TextValue "D50" "22.58"
$ws.Range("E50").Value = "  -2.49%  "
$ws.Range("E51").Value = "  +0.13%  "

$excel.CutCopyMode = $false

